# Fix Barclay XLS format
#
# The "Umsaetze_nok_wrongdate1.xlsx" sample sheet used two overly verbose /
# inconsistent column headings in the Barclaycard transaction export, and the
# single data row stored a stray literal ("Haendler") instead of reusing the
# same text as its column header ("Haendlerdetails"). Align the wording with
# the current Barclaycard export format:
#
#   - "Name des Karteninhabers" -> "Karteninhaber"
#   - "Haendlerdetails"         -> "Details"
#   - data cell under "Details" now reuses that same "Details" text instead
#     of the old standalone "Haendler" label.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row 13: column L ("Name des Karteninhabers" -> "Karteninhaber")
$ws.Cells.Item(13, 12).Value = "Karteninhaber"

# Header row 13: column O ("Haendlerdetails" -> "Details")
$ws.Cells.Item(13, 15).Value = "Details"

# Data row 14: column O used to contain the standalone word "Haendler";
# it now matches the (renamed) header text above it, "Details".
$ws.Cells.Item(14, 15).Value = "Details"

Write-Host "Barclaycard header labels updated (Karteninhaber / Details)"
